$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: swap home/away match data between mis-ordered row pairs ---
# Each pair below has its F:V block (match details) fully swapped; A (index) and E (date) stay put.
$tmp1 = $ws.Range("F27:V27").Value()
$tmp2 = $ws.Range("F28:V28").Value()
$ws.Range("F27:V27").Value = $tmp2
$ws.Range("F28:V28").Value = $tmp1

$tmp1 = $ws.Range("F38:V38").Value()
$tmp2 = $ws.Range("F39:V39").Value()
$ws.Range("F38:V38").Value = $tmp2
$ws.Range("F39:V39").Value = $tmp1

$tmp1 = $ws.Range("F47:V47").Value()
$tmp2 = $ws.Range("F48:V48").Value()
$ws.Range("F47:V47").Value = $tmp2
$ws.Range("F48:V48").Value = $tmp1

$tmp1 = $ws.Range("F66:V66").Value()
$tmp2 = $ws.Range("F67:V67").Value()
$ws.Range("F66:V66").Value = $tmp2
$ws.Range("F67:V67").Value = $tmp1

$tmp1 = $ws.Range("F87:V87").Value()
$tmp2 = $ws.Range("F88:V88").Value()
$ws.Range("F87:V87").Value = $tmp2
$ws.Range("F88:V88").Value = $tmp1

$tmp1 = $ws.Range("F103:V103").Value()
$tmp2 = $ws.Range("F104:V104").Value()
$ws.Range("F103:V103").Value = $tmp2
$ws.Range("F104:V104").Value = $tmp1

$tmp1 = $ws.Range("F115:V115").Value()
$tmp2 = $ws.Range("F116:V116").Value()
$ws.Range("F115:V115").Value = $tmp2
$ws.Range("F116:V116").Value = $tmp1

# --- Part 2: append 7 new match rows (175-181) at the end of the table ---
# Seed each new row from the last existing row (174) via Copy so number formats / styles
# (bordered index column, date-time column) and the text-typed "temporada" cell carry over
# exactly (column D holds the text "2023", identical to row 174, so it never needs rewriting),
# then overwrite the per-row content.

for ($r = 175; $r -le 181; $r++) {
    $ws.Range("A174:V174").Copy($ws.Range("A${r}:V${r}"))
}

# Row 175
$ws.Range("A175").Value = 174
$ws.Range("E175").Value = 45228.41666666666
$ws.Range("F175").Value = "Aksu"
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = "Tobol"
$ws.Range("I175").Value = 3
$ws.Range("J175").Value = 2.62
$ws.Range("K175").Value = "27/10/2023 23:12"
$ws.Range("L175").Value = 2.51
$ws.Range("M175").Value = "29/10/2023 09:57"
$ws.Range("N175").Value = 2.95
$ws.Range("O175").Value = "27/10/2023 23:12"
$ws.Range("P175").Value = 3.55
$ws.Range("Q175").Value = "29/10/2023 09:54"
$ws.Range("R175").Value = 2.34
$ws.Range("S175").Value = "27/10/2023 23:12"
$ws.Range("T175").Value = 2.35
$ws.Range("U175").Value = "29/10/2023 09:57"
$ws.Range("V175").Value = "https://www.betexplorer.com/football/kazakhstan/premier-league/aksu-tobol/CUOV8gZn/"

# Row 176
$ws.Range("A176").Value = 175
$ws.Range("E176").Value = 45228.41666666666
$ws.Range("F176").Value = "Atyrau"
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = "Kairat Almaty"
$ws.Range("I176").Value = 0
$ws.Range("J176").Value = 2.92
$ws.Range("K176").Value = "27/10/2023 23:12"
$ws.Range("L176").Value = 2.63
$ws.Range("M176").Value = "29/10/2023 09:53"
$ws.Range("N176").Value = 2.95
$ws.Range("O176").Value = "27/10/2023 23:12"
$ws.Range("P176").Value = 3.16
$ws.Range("Q176").Value = "29/10/2023 09:53"
$ws.Range("R176").Value = 2.15
$ws.Range("S176").Value = "27/10/2023 23:12"
$ws.Range("T176").Value = 2.45
$ws.Range("U176").Value = "29/10/2023 09:53"
$ws.Range("V176").Value = "https://www.betexplorer.com/football/kazakhstan/premier-league/atyrau-kairat-almaty/I7PR9ZJt/"

# Row 177
$ws.Range("A177").Value = 176
$ws.Range("E177").Value = 45228.41666666666
$ws.Range("F177").Value = "Kyzylzhar"
$ws.Range("G177").Value = 0
$ws.Range("H177").Value = "Kaisar Kyzylorda"
$ws.Range("I177").Value = 1
$ws.Range("J177").Value = 1.79
$ws.Range("K177").Value = "27/10/2023 23:12"
$ws.Range("L177").Value = 1.81
$ws.Range("M177").Value = "29/10/2023 09:51"
$ws.Range("N177").Value = 3.02
$ws.Range("O177").Value = "27/10/2023 23:12"
$ws.Range("P177").Value = 3.26
$ws.Range("Q177").Value = "29/10/2023 09:52"
$ws.Range("R177").Value = 3.87
$ws.Range("S177").Value = "27/10/2023 23:12"
$ws.Range("T177").Value = 4.04
$ws.Range("U177").Value = "29/10/2023 09:48"
$ws.Range("V177").Value = "https://www.betexplorer.com/football/kazakhstan/premier-league/kyzylzhar-kaisar-kyzylorda/Gbw0EFBP/"

# Row 178
$ws.Range("A178").Value = 177
$ws.Range("E178").Value = 45228.5
$ws.Range("F178").Value = "Okzhetpes"
$ws.Range("G178").Value = 1
$ws.Range("H178").Value = "Maqtaaral"
$ws.Range("I178").Value = 1
$ws.Range("J178").Value = 2.58
$ws.Range("K178").Value = "28/10/2023 01:42"
$ws.Range("L178").Value = 2.35
$ws.Range("M178").Value = "29/10/2023 11:52"
$ws.Range("N178").Value = 2.91
$ws.Range("O178").Value = "28/10/2023 01:42"
$ws.Range("P178").Value = 3.03
$ws.Range("Q178").Value = "29/10/2023 11:52"
$ws.Range("R178").Value = 2.41
$ws.Range("S178").Value = "28/10/2023 01:42"
$ws.Range("T178").Value = 2.86
$ws.Range("U178").Value = "29/10/2023 11:52"
$ws.Range("V178").Value = "https://www.betexplorer.com/football/kazakhstan/premier-league/okzhetpes-maqtaaral/G6XOkWYO/"

# Row 179
$ws.Range("A179").Value = 178
$ws.Range("E179").Value = 45228.5
$ws.Range("F179").Value = "Aktobe"
$ws.Range("G179").Value = 2
$ws.Range("H179").Value = "Zhetysu Taldykorgan"
$ws.Range("I179").Value = 2
$ws.Range("J179").Value = 1.3
$ws.Range("K179").Value = "28/10/2023 01:42"
$ws.Range("L179").Value = 1.44
$ws.Range("M179").Value = "29/10/2023 11:54"
$ws.Range("N179").Value = 4.38
$ws.Range("O179").Value = "28/10/2023 01:42"
$ws.Range("P179").Value = 4.46
$ws.Range("Q179").Value = "29/10/2023 11:56"
$ws.Range("R179").Value = 6.48
$ws.Range("S179").Value = "28/10/2023 01:42"
$ws.Range("T179").Value = 5.31
$ws.Range("U179").Value = "29/10/2023 11:56"
$ws.Range("V179").Value = "https://www.betexplorer.com/football/kazakhstan/premier-league/aktobe-zhetysu-taldykorgan/nmHr6iJ5/"

# Row 180
$ws.Range("A180").Value = 179
$ws.Range("E180").Value = 45228.5
$ws.Range("F180").Value = "Kaspij Aktau"
$ws.Range("G180").Value = 1
$ws.Range("H180").Value = "Ordabasy"
$ws.Range("I180").Value = 0
$ws.Range("J180").Value = 5.31
$ws.Range("K180").Value = "28/10/2023 01:42"
$ws.Range("L180").Value = 4.08
$ws.Range("M180").Value = "29/10/2023 11:54"
$ws.Range("N180").Value = 3.92
$ws.Range("O180").Value = "28/10/2023 01:42"
$ws.Range("P180").Value = 3.99
$ws.Range("Q180").Value = "29/10/2023 11:54"
$ws.Range("R180").Value = 1.41
$ws.Range("S180").Value = "28/10/2023 01:42"
$ws.Range("T180").Value = 1.64
$ws.Range("U180").Value = "29/10/2023 11:54"
$ws.Range("V180").Value = "https://www.betexplorer.com/football/kazakhstan/premier-league/kaspij-aktau-ordabasy/QXSZ7Dlg/"

# Row 181
$ws.Range("A181").Value = 180
$ws.Range("E181").Value = 45228.5
$ws.Range("F181").Value = "Shakhter Karagandy"
$ws.Range("G181").Value = 1
$ws.Range("H181").Value = "FC Astana"
$ws.Range("I181").Value = 1
$ws.Range("J181").Value = 3.24
$ws.Range("K181").Value = "29/10/2023 10:10"
$ws.Range("L181").Value = 3.24
$ws.Range("M181").Value = "29/10/2023 10:10"
$ws.Range("N181").Value = 2.13
$ws.Range("O181").Value = "29/10/2023 10:10"
$ws.Range("P181").Value = 2.13
$ws.Range("Q181").Value = "29/10/2023 10:10"
$ws.Range("R181").Value = 2.83
$ws.Range("S181").Value = "29/10/2023 10:10"
$ws.Range("T181").Value = 2.83
$ws.Range("U181").Value = "29/10/2023 10:10"
$ws.Range("V181").Value = "https://www.betexplorer.com/football/kazakhstan/premier-league/shakhter-karagandy-fc-astana/8dIv7X3a/"
